$d = $word.ActiveDocument

# --- 1) First paragraph: append two trailing spaces to the existing text,
#        then append a red-colored annotation split across three runs.
$firstPara = $d.Paragraphs(1).Range
# Exclude the trailing paragraph mark from the range so the insertion lands
# at the very end of the paragraph's text, not the start of the next one.
$firstText = $d.Range($firstPara.Start, $firstPara.End - 1)
$firstText.Collapse(0)

$firstText.InsertAfter("  ")
$firstText.Collapse(0)

$firstText.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$firstText.Font.Color = 255
$firstText.Collapse(0)

$firstText.InsertAfter("rsion for main branch")
$firstText.Font.Color = 255
$firstText.Collapse(0)

$firstText.InsertAfter(")")
$firstText.Font.Color = 255
$firstText.Collapse(0)

# --- 2) Remove the trailing "ank God almighty, we are free at last."
#        paragraph entirely (it was the last paragraph in the document).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$delRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)
$delRange.Delete()
